# Adds a "minimum utilization rate" variable to the hubs sheet and a
# "goal of co2 emissions" variable to the emissions sheet, and switches
# the active sheet from "cost" to "hubs".

$wb = $excel.ActiveWorkbook

$emissions = $wb.Worksheets.Item("emissions")
$emissions.Range("A3").Value = "goal of co2 emissions"
$emissions.Range("B3").Value = 1000000
$emissions.Range("C3").Value = "CO2"

$hubs = $wb.Worksheets.Item("hubs")
$hubs.Range("A3").Value = "minimum utilization rate"
$hubs.Range("B3").Value = 0.6
$hubs.Range("C3").Value = "Dmnl"

$emissions.Activate()
$emissions.Range("C4").Select()

$hubs.Activate()
$hubs.Range("C3").Select()
